# This script rewrites the "Analysis Results" sheet to match the updated
# report layout:
#   - A new "Date and Time" row is inserted at the very top (row 1), pushing
#     every existing row down by one.
#   - "Highest Cell Voltage" / "Lowest Cell Voltage" are swapped in order,
#     and likewise "highest cell temp" / "lowest cell temp".
#   - Several row labels gain unit suffixes, e.g. "(kW)", "(C)", "(V)", "-BMS".
#   - The "Maximum BMS Temperature in C" row is removed.
#   - A new "Cycle Count of battery" row is inserted before "Idling time percentage".
#   - Two new rows, "Time spent in 70-80 km/h" and "Time spent in 80-90 km/h",
#     are appended at the bottom.
#   - A handful of numeric values were recalculated/corrected.
#
# Because almost every row label and/or value changes position or content,
# the simplest reliable approach is to clear the sheet and rewrite every
# cell explicitly in its final location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Clear()

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 15:26:19.887000 to 2024-03-12 17:37:35.929000"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = [double]"0.09123127314814815"
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = [double]"33.47750777777778"

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = [double]"1736.036073496667"

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = [double]"37.326"

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = [double]"5.944"

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = [double]"100"

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = [double]"15"

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = [double]"70.59788385330441"

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = [double]"24.59048315249763"

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = [double]"85"

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Eco mode`n77.10%`nSports mode`n18.20%`nCustom mode`n0.05%"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = [double]"5317.71037"

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = [double]"-799.4026432064467"

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = [double]"69.98707205138889"

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = [double]"3.875203494701095"

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = [double]"3.522"

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = [double]"3.05"

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = [double]"0.472"

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = [double]"38"

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = [double]"48"

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = [double]"10"

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = [double]"63"

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = [double]"67"

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = [double]"61"

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = [double]"46"

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = [double]"0"

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = [double]"0"

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = [double]"48"

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = [double]"38"

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = [double]"10"

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = [double]"56"

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = [double]"1.874740435555556"

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = [double]"6.606968181917856e-08"

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = [double]"137"

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = [double]"4.941536748329622"

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = [double]"7.743596881959911"

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = [double]"6.531180400890868"

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = [double]"8.083240534521158"

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = [double]"48.46185968819599"

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = [double]"22.46380846325167"

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = [double]"0"

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = [double]"0"

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = [double]"0"

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = [double]"0"
